{"js": "// Corre\u00e7\u00e3o ate o DFD ac5\n// Remove the two empty paragraphs that sit between the title\n// (\"Lista de Caracter\u00edsticas\") and the \"Legenda:\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the title paragraph and the \"Legenda:\" paragraph, then delete every\n// empty paragraph strictly between them (there are two in this document).\nconst items = paragraphs.items;\nlet titleIndex = -1;\nlet legendaIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (titleIndex === -1 && text.indexOf(\"Lista de Caracter\u00edsticas\") !== -1) {\n    titleIndex = i;\n  } else if (titleIndex !== -1 && text.indexOf(\"Legenda:\") !== -1) {\n    legendaIndex = i;\n    break;\n  }\n}\n\nif (titleIndex !== -1 && legendaIndex !== -1) {\n  for (let i = legendaIndex - 1; i > titleIndex; i--) {\n    if (items[i].text === \"\") {\n      items[i].delete();\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Corre\u00e7\u00e3o ate o DFD ac5\n# Remove the two empty paragraphs that sit between the title\n# (\"Lista de Caracter\u00edsticas\") and the \"Legenda:\" paragraph.\n\n$d = $word.ActiveDocument\n\n$titleIndex = -1\n$legendaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($titleIndex -eq -1 -and $t -like \"*Lista de Caracter*\") {\n        $titleIndex = $i\n    } elseif ($titleIndex -ne -1 -and $t -like \"*Legenda:*\") {\n        $legendaIndex = $i\n        break\n    }\n}\n\nif ($titleIndex -ne -1 -and $legendaIndex -ne -1) {\n    for ($i = $legendaIndex - 1; $i -gt $titleIndex; $i--) {\n        $p = $d.Paragraphs($i)\n        $txt = $p.Range.Text.Trim([char]13, [char]7)\n        if ($txt -eq \"\") {\n            $p.Range.Delete()\n        }\n    }\n}\n"}
